$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F ("想去人数") values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3095
$ws1.Range("F3").Value = 503
$ws1.Range("F4").Value = 68
$ws1.Range("F5").Value = 59
$ws1.Range("F6").Value = 15
$ws1.Range("F9").Value = 1078
$ws1.Range("F10").Value = 15134
$ws1.Range("F12").Value = 150
$ws1.Range("F14").Value = 6000
$ws1.Range("F15").Value = 611
$ws1.Range("F17").Value = 55
$ws1.Range("F19").Value = 1251
$ws1.Range("F21").Value = 104
$ws1.Range("F22").Value = 6
$ws1.Range("F23").Value = 204
$ws1.Range("F24").Value = 834
$ws1.Range("F25").Value = 2992
$ws1.Range("F27").Value = 10841
$ws1.Range("F29").Value = 2
$ws1.Range("F31").Value = 140
$ws1.Range("F32").Value = 3770
$ws1.Range("F33").Value = 256

# Sheet "全部类型" (sheet4) - update column F ("想去人数") values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3095
$ws4.Range("F4").Value = 503
$ws4.Range("F5").Value = 68
$ws4.Range("F6").Value = 59
$ws4.Range("F7").Value = 15
$ws4.Range("F10").Value = 1078
$ws4.Range("F11").Value = 15135
$ws4.Range("F13").Value = 150
$ws4.Range("F15").Value = 6000
$ws4.Range("F16").Value = 611
$ws4.Range("F18").Value = 55
$ws4.Range("F20").Value = 1251
$ws4.Range("F22").Value = 104
$ws4.Range("F23").Value = 6
$ws4.Range("F24").Value = 204
$ws4.Range("F25").Value = 834
$ws4.Range("F26").Value = 2992
$ws4.Range("F29").Value = 10841
$ws4.Range("F31").Value = 2
$ws4.Range("F33").Value = 140
$ws4.Range("F34").Value = 3770
$ws4.Range("F35").Value = 256

Write-Output "Updated F column values on 展览 and 全部类型 sheets"
